$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51
$ws.Range("A51").Value = 49
$ws.Range("B51").Value = 10
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0.0003
$ws.Range("H51").Value = "Regular"
$ws.Range("L51").Value = "<function relu at 0x11d707488>"
$ws.Range("M51").Value = 0.9351999759674072
$ws.Range("N51").Value = 0.3497999906539917
$ws.Range("P51").Value = 0.2422611862421036
$ws.Range("Q51").Value = 3.15626335144043

# Row 52
$ws.Range("A52").Value = 50
$ws.Range("B52").Value = 30
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("E52").Value = 0
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 0.0003
$ws.Range("H52").Value = "Regular"
$ws.Range("L52").Value = "<function relu at 0x11d707488>"
$ws.Range("M52").Value = 0.9556999802589417
$ws.Range("N52").Value = 0.04399999976158142
$ws.Range("P52").Value = 0.2216933816671371
$ws.Range("Q52").Value = 23.85161018371582
$ws.Range("T52").Value = "weights/model_367.ckpt"
